$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) next to the existing "sum" column (G).
# Copy the header cell's formatting (bold font, border, centered alignment)
# so the new header reuses the same style as the other headers.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the data value for row 2 under the new column.
$ws.Range("H2").Value = 1
